# Fix syllabus for cm016
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 corresponds to cm016. Fix the typo "tree based" -> "tree-based"
# and flip the link_it flag from False to True.
$ws.Range("D17").Value = "Statistical learning: resampling and tree-based methods"
$ws.Range("C17").Value = $true

# Update the active selection shown in the sheet view.
$ws.Range("D18").Select()
